$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1195
$ws1.Range("F7").Value = 4391
$ws1.Range("F8").Value = 2621
$ws1.Range("F10").Value = 2542
$ws1.Range("F15").Value = 667
$ws1.Range("F18").Value = 333
$ws1.Range("F26").Value = 563
$ws1.Range("F27").Value = 701
$ws1.Range("F28").Value = 110
$ws1.Range("F31").Value = 1624
$ws1.Range("F32").Value = 1076
$ws1.Range("F33").Value = 152
$ws1.Range("F35").Value = 1168
$ws1.Range("F36").Value = 2082
$ws1.Range("F37").Value = 279
$ws1.Range("F39").Value = 552
$ws1.Range("F41").Value = 28
$ws1.Range("F43").Value = 670
$ws1.Range("F44").Value = 1342
$ws1.Range("F45").Value = 113
$ws1.Range("F47").Value = 443
$ws1.Range("F48").Value = 73

# Sheet "演出" (sheet2): update F column values
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 11
$ws2.Range("F11").Value = 10

# Sheet "全部类型" (sheet4): update F column values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1195
$ws4.Range("F6").Value = 4391
$ws4.Range("F7").Value = 2621
$ws4.Range("F8").Value = 2542
$ws4.Range("F12").Value = 667
$ws4.Range("F15").Value = 333
$ws4.Range("F22").Value = 563
$ws4.Range("F23").Value = 701
$ws4.Range("F24").Value = 110
$ws4.Range("F29").Value = 1624
$ws4.Range("F30").Value = 1076
$ws4.Range("F31").Value = 152
$ws4.Range("F34").Value = 2082
$ws4.Range("F35").Value = 279
$ws4.Range("F37").Value = 11
$ws4.Range("F39").Value = 10
$ws4.Range("F40").Value = 552
$ws4.Range("F42").Value = 28
$ws4.Range("F43").Value = 670
$ws4.Range("F44").Value = 1342
$ws4.Range("F46").Value = 113
$ws4.Range("F47").Value = 443
$ws4.Range("F48").Value = 73

$wb.Save()
